$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ C = val; D = val; E = val } (only include cells that change)
$updates = @{
    3  = @{ C = 2;  E = 0.04 }
    4  = @{ C = 1;  D = 1;  E = 0.37 }
    5  = @{ C = 5;  D = 18 }
    6  = @{ C = 2;  D = 73 }
    7  = @{ D = 4;  E = 0.09 }
    11 = @{ C = 41; D = 25; E = 0 }
    13 = @{ C = 0;  D = 0;  E = 1 }
    14 = @{ C = 8;  D = 3;  E = 0.03 }
    19 = @{ C = 8;  D = 9;  E = 0.12 }
    20 = @{ C = 1;  E = 0.37 }
    21 = @{ C = 8;  E = 0.06 }
    24 = @{ C = 2;  D = 1;  E = 0.27 }
    25 = @{ D = 0;  E = 1 }
    26 = @{ D = 0;  E = 1 }
    27 = @{ C = 4;  D = 4;  E = 0.2 }
    29 = @{ D = 2;  E = 0 }
    31 = @{ C = 1;  D = 1;  E = 0.37 }
    33 = @{ D = 0;  E = 1 }
    34 = @{ D = 0;  E = 1 }
    35 = @{ C = 2;  E = 0.27 }
    36 = @{ C = 6;  D = 3 }
    37 = @{ C = 7;  D = 1;  E = 0.01 }
    38 = @{ C = 9;  D = 4;  E = 0.03 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cellAddr = "$col$row"
        $ws.Range($cellAddr).Value = $cols[$col]
    }
}
